$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the new data row (row 2: B2=2, C2=3, D2=4)
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 4

# Apply AutoFilter over A1:J12 (header row 1, data rows 2:12)
$ws.Range("A1:J12").AutoFilter()
